# uGrid_Input.xlsx edit: "Additional of solar calcs and RRT attempt"
#
# Summary of the change being applied:
#  1. PSO sheet: rename the repeatability-test flag stored in O2 from
#     "test1" to "repeatabilityTest3".
#  2. Econ sheet: add a new "tariff_hillclimb_multiplier" parameter in
#     column AK.
#  3. Tech sheet: add a new "trans_losses" parameter in column I.
#  4. Add a brand-new "Solar" sheet (after "Tech") holding the solar/PV
#     calculation inputs (year, longitude, latitude, timezone, slope,
#     azimuth, pg, fpv, alpha_p, eff_mpp, f_inv).
#  5. Leave "Tech" as the active/selected sheet, matching the saved file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Econ!AK1:AK2 - new tariff hillclimb multiplier column
# ---------------------------------------------------------------------
$wsEcon = $wb.Worksheets.Item("Econ")
$wsEcon.Range("AK1").Value = "tariff_hillclimb_multiplier"
$wsEcon.Range("AK2").Value = 1.01

# ---------------------------------------------------------------------
# 2) New "Solar" sheet, inserted right after "Tech"
# ---------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("Tech")
$wsSolar = $wb.Worksheets.Add([Type]::Missing, $wsTech)
$wsSolar.Name = "Solar"

# Headers typed in this order: B, C, A, D, E, F, G, H, (later) J, I, (later) K
$wsSolar.Range("B1").Value = "longitude"
$wsSolar.Range("C1").Value = "latitude"
$wsSolar.Range("A1").Value = "year"
$wsSolar.Range("D1").Value = "timezone"
$wsSolar.Range("E1").Value = "slope"
$wsSolar.Range("F1").Value = "azimuth"
$wsSolar.Range("G1").Value = "pg"
$wsSolar.Range("H1").Value = "fpv"
$wsSolar.Range("I1").Value = "alpha_p"
$wsSolar.Range("J1").Value = "eff_mpp"

# Row 2 data values
$wsSolar.Range("A2").Value = 2005
$wsSolar.Range("B2").Value = -33
$wsSolar.Range("C2").Value = 18
$wsSolar.Range("D2").Value = 2
$wsSolar.Range("E2").Value = 0
$wsSolar.Range("F2").Value = 0
$wsSolar.Range("G2").Value = 0.2
$wsSolar.Range("H2").Value = 0.9
$wsSolar.Range("I2").Value = -0.002
$wsSolar.Range("J2").Value = 0.9

# ---------------------------------------------------------------------
# 3) PSO!O2 - rename the old "test1" flag to "repeatabilityTest3"
# ---------------------------------------------------------------------
$wsPSO = $wb.Worksheets.Item("PSO")
$wsPSO.Range("O2").Value = "repeatabilityTest3"

# ---------------------------------------------------------------------
# 4) Solar!K1:K2 - f_inv, added after the PSO edit above
# ---------------------------------------------------------------------
$wsSolar.Range("K1").Value = "f_inv"
$wsSolar.Range("K2").Value = 0.9

# ---------------------------------------------------------------------
# 5) Tech!I1:I2 - new trans_losses column
# ---------------------------------------------------------------------
$wsTech.Range("I1").Value = "trans_losses"
$wsTech.Range("I2").Value = 0.08

# ---------------------------------------------------------------------
# 6) Column widths on PSO - autofit to content (matches the bestFit
#    column metadata recorded in the saved workbook).
# ---------------------------------------------------------------------
[void]$wsPSO.Range("A1:O2").EntireColumn.AutoFit()

# ---------------------------------------------------------------------
# 7) Restore each sheet's on-screen selection to match the saved file,
#    then leave "Tech" as the active sheet/tab.
# ---------------------------------------------------------------------
[void]$wsPSO.Range("N14").Select()
[void]$wsEcon.Range("AE9").Select()
[void]$wsSolar.Range("K2").Select()
[void]$wsTech.Range("I3").Select()
[void]$wsTech.Activate()
